# Update "Contribution Summary - Final.xlsx" Sheet1 with the revised
# contribution breakdown: more granular tasks + a new "Contribution
# Details" column (C) filled in for Owen Randolph's rows, and three
# additional rows appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table (rows 2-16), columns: Name, Task, Contribution Details, Hours
$data = @(
    @(2,  "Owen Randolph",    "Develop Cloud Infrastructure",           "Set up and configure AWS Infrastructure: IAM, VPC, security groups, subnets, availability zones, Aurora MySQL, S3", 4),
    @(3,  "Owen Randolph",    "Connect MySQL Workbench db to Aurora",   "Use routing through MySQL workbench to run the database on Aurora", 1),
    @(4,  "Owen Randolph",    "Connect S3 to and configure Quicksight", "Data ingestion for use in Quicksight, set up normalized tables in Quicksight for dashbord use", 4),
    @(5,  "Owen Randolph",    "Report Draft",                           "Write notes on project build, add screenshots, edit and resize pictures and text", 3),
    @(6,  "Owen Randolph",    "Report Detail- Technical Description",   "refine formatting and clarity of notes on technical build part of the project", 3),
    @(7,  "Owen Randolph",    "Report Detail- Future Expansion",        "data engineering features and cloud architecture services to offer expansion plans", 0.5),
    @(8,  "Owen Randolph",    "Report Detail- References",              "Used hyperlinks", 0.5),
    @(9,  "Owen Randolph",    "Powerpoint Creation",                    "Added screenshots, wrote and formatted text, added notes for presentation", 2),
    @(10, "Marcos Fernandez", "Region Analysis Dashboard Development",  $null, 2.5),
    @(11, "Marcos Fernandez", "Report Detail- App Functionality",       $null, 1.5),
    @(12, "Marcos Fernandez", "Report Detail - Reflections",            $null, 1.5),
    @(13, "Marcos Fernandez", "Report Finalize",                        $null, 1),
    @(14, "Gabriel Tharp",    "Product Analysis Dashboard Development", $null, 2.5),
    @(15, "Gabriel Tharp",    "Report Detail- Purpose & Audience ",     $null, 1.5),
    @(16, "Gabriel Tharp",    "Report Detail - SQL Queries",            $null, 1.5)
)

foreach ($row in $data) {
    $r       = $row[0]
    $name    = $row[1]
    $task    = $row[2]
    $details = $row[3]
    $hours   = $row[4]

    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $task

    if ($details) {
        $ws.Cells.Item($r, 3).Value = $details
    } else {
        $ws.Cells.Item($r, 3).Value = ""
    }

    $ws.Cells.Item($r, 4).Value = $hours
}
